$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = "2026-03-01 07:48:26"
$ws.Range("H2").Value = "'91%"
$ws.Range("N2").Value = "-2.6 °C 7:15 TU"
$ws.Range("O2").Value = "-1.3 °C"
$ws.Range("E3").Value = "2026-03-01 07:48:28"
$ws.Range("L3").Value = "21.2 km/h - 110º 7:03 TU"
$ws.Range("N3").Value = "-4.2 °C 7:01 TU"
$ws.Range("O3").Value = "-3.8 °C"
$ws.Range("E4").Value = "2026-03-01 07:48:30"
$ws.Range("J4").Value = "1025.7 hPa"
$ws.Range("M4").Value = "8.9 °C 7:29 TU"
$ws.Range("E5").Value = "2026-03-01 07:48:33"
$ws.Range("E6").Value = "2026-03-01 07:48:35"
$ws.Range("K6").Value = "0.1 MJ/m2"
$ws.Range("O6").Value = "9.4 °C"
$ws.Range("E7").Value = "2026-03-01 07:48:38"
$ws.Range("K7").Value = "0.1 MJ/m2"
$ws.Range("E8").Value = "2026-03-01 07:48:40"
$ws.Range("K8").Value = "0.1 MJ/m2"
$ws.Range("E9").Value = "2026-03-01 07:48:42"
$ws.Range("H9").Value = "'62%"
$ws.Range("K9").Value = "0.1 MJ/m2"
$ws.Range("N9").Value = "10.3 °C 7:28 TU"
$ws.Range("O9").Value = "11.5 °C"
$ws.Range("E10").Value = "2026-03-01 07:48:45"
$ws.Range("K10").Value = "0.1 MJ/m2"
$ws.Range("E11").Value = "2026-03-01 07:48:47"
$ws.Range("N11").Value = "5.8 °C 7:00 TU"
$ws.Range("E12").Value = "2026-03-01 07:48:50"
$ws.Range("H12").Value = "'78%"
$ws.Range("E13").Value = "2026-03-01 07:48:52"
$ws.Range("J13").Value = "1026.4 hPa"
$ws.Range("K13").Value = "0.1 MJ/m2"
$ws.Range("O13").Value = "4.0 °C"
$ws.Range("E14").Value = "2026-03-01 07:48:55"
$ws.Range("K14").Value = "0.1 MJ/m2"
$ws.Range("E15").Value = "2026-03-01 07:48:57"
$ws.Range("O15").Value = "8.7 °C"
$ws.Range("E16").Value = "2026-03-01 07:48:59"
$ws.Range("H16").Value = "'82%"
$ws.Range("E17").Value = "2026-03-01 07:49:02"
$ws.Range("E18").Value = "2026-03-01 07:49:04"
$ws.Range("K18").Value = "0.1 MJ/m2"
$ws.Range("M18").Value = "8.4 °C 7:29 TU"
$ws.Range("O18").Value = "6.9 °C"
$ws.Range("E19").Value = "2026-03-01 07:49:07"
$ws.Range("I19").Value = "1.1 mm"
$ws.Range("N19").Value = "5.8 °C 7:00 TU"
$ws.Range("E20").Value = "2026-03-01 07:49:09"
$ws.Range("K20").Value = "0.1 MJ/m2"
$ws.Range("O20").Value = "-3.4 °C"
$ws.Range("E21").Value = "2026-03-01 07:49:12"
$ws.Range("H21").Value = "'90%"
$ws.Range("K21").Value = "0.1 MJ/m2"
$ws.Range("L21").Value = "7.2 km/h - 216º 7:25 TU"
$ws.Range("O21").Value = "6.1 °C"
$ws.Range("E22").Value = "2026-03-01 07:49:14"
$ws.Range("K22").Value = "0.2 MJ/m2"
$ws.Range("O22").Value = "-5.6 °C"
$ws.Range("E23").Value = "2026-03-01 07:49:17"
$ws.Range("K23").Value = "0.1 MJ/m2"
$ws.Range("N23").Value = "-4.9 °C 7:25 TU"
$ws.Range("E24").Value = "2026-03-01 07:49:19"
$ws.Range("M24").Value = "6.5 °C 7:29 TU"
$ws.Range("O24").Value = "4.9 °C"
$ws.Range("E25").Value = "2026-03-01 07:49:21"
$ws.Range("H25").Value = "'91%"
$ws.Range("K25").Value = "0.3 MJ/m2"
$ws.Range("E26").Value = "2026-03-01 07:49:24"
$ws.Range("E27").Value = "2026-03-01 07:49:26"
$ws.Range("K27").Value = "0.2 MJ/m2"
$ws.Range("L27").Value = "10.4 km/h - 205º 7:27 TU"
$ws.Range("N27").Value = "-2.9 °C 7:10 TU"
$ws.Range("O27").Value = "-1.7 °C"
$ws.Range("E28").Value = "2026-03-01 07:49:29"
$ws.Range("O28").Value = "8.4 °C"
$ws.Range("E29").Value = "2026-03-01 07:49:31"
$ws.Range("H29").Value = "'82%"
$ws.Range("K29").Value = "0.1 MJ/m2"
$ws.Range("O29").Value = "8.6 °C"
$ws.Range("E30").Value = "2026-03-01 07:49:34"
$ws.Range("H30").Value = "'81%"
$ws.Range("J30").Value = "1025.6 hPa"
$ws.Range("K30").Value = "0.1 MJ/m2"
$ws.Range("O30").Value = "10.0 °C"
$ws.Range("E31").Value = "2026-03-01 07:49:37"
$ws.Range("J31").Value = "1024.7 hPa"
$ws.Range("E32").Value = "2026-03-01 07:49:39"
$ws.Range("M32").Value = "5.2 °C 7:15 TU"
$ws.Range("O32").Value = "3.2 °C"
$ws.Range("E33").Value = "2026-03-01 07:49:42"
$ws.Range("J33").Value = "1025.7 hPa"
$ws.Range("K33").Value = "0.1 MJ/m2"
$ws.Range("N33").Value = "3.5 °C 7:17 TU"
$ws.Range("E34").Value = "2026-03-01 07:49:44"
$ws.Range("K34").Value = "0.1 MJ/m2"
$ws.Range("N34").Value = "-1.0 °C 7:23 TU"
$ws.Range("E35").Value = "2026-03-01 07:49:47"
$ws.Range("E36").Value = "2026-03-01 07:49:50"
$ws.Range("H36").Value = "'72%"
$ws.Range("K36").Value = "0.2 MJ/m2"
$ws.Range("L36").Value = "22.0 km/h - 26º 7:01 TU"
$ws.Range("M36").Value = "13.1 °C 7:24 TU"
$ws.Range("O36").Value = "10.6 °C"
$ws.Range("E37").Value = "2026-03-01 07:49:52"
$ws.Range("O37").Value = "6.3 °C"
$ws.Range("E38").Value = "2026-03-01 07:49:55"
$ws.Range("K38").Value = "0.1 MJ/m2"
$ws.Range("M38").Value = "9.7 °C 7:11 TU"
$ws.Range("O38").Value = "8.8 °C"
$ws.Range("E39").Value = "2026-03-01 07:49:58"
$ws.Range("K39").Value = "0.1 MJ/m2"
$ws.Range("N39").Value = "-4.7 °C 7:29 TU"
$ws.Range("O39").Value = "-3.3 °C"
$ws.Range("E40").Value = "2026-03-01 07:50:00"
$ws.Range("J40").Value = "1025.9 hPa"
$ws.Range("O40").Value = "6.7 °C"
$ws.Range("E41").Value = "2026-03-01 07:50:03"
$ws.Range("E42").Value = "2026-03-01 07:50:05"
$ws.Range("H42").Value = "'90%"
$ws.Range("O42").Value = "7.8 °C"
$ws.Range("E43").Value = "2026-03-01 07:50:08"
$ws.Range("N43").Value = "8.2 °C 7:00 TU"
$ws.Range("E44").Value = "2026-03-01 07:50:10"
$ws.Range("K44").Value = "0.1 MJ/m2"
$ws.Range("N44").Value = "-4.5 °C 7:29 TU"
$ws.Range("O44").Value = "-3.1 °C"
$ws.Range("E45").Value = "2026-03-01 07:50:13"
$ws.Range("J45").Value = "1027.2 hPa"
$ws.Range("O45").Value = "3.0 °C"
$ws.Range("E46").Value = "2026-03-01 07:50:16"
$ws.Range("M46").Value = "8.7 °C 7:29 TU"
